$p = $ppt.ActivePresentation

# --- Slide 2: "Task:" paragraph placeholder ------------------------------
# Fix the "calssification on a" -> "classification on a" typo and drop the
# stray leading space that preceded "on a " once the previous run gains its
# own trailing space.
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$fullText = $tr.Text
$typo = "calssification"
$idx = $fullText.IndexOf($typo)
if ($idx -ge 0) {
    $startPos = $idx + 1

    # Remove the misspelled run entirely; it merges into the following
    # " on a " run (which carries the plain, non-error formatting).
    $typoRun = $tr.Characters($startPos, $typo.Length)
    $typoRun.Text = ""

    # Re-type the correct word (with its trailing space) over the single
    # leftover space character; this produces a new run using the
    # surrounding (error-free) run formatting.
    $spaceChar = $tr.Characters($startPos, 1)
    $spaceChar.Text = "classification "

    # The following run currently reads " on a " (leading space); trim it
    # down to "on a " since "classification " now supplies its own space.
    $afterLen = ("on a ").Length + 1
    $afterRun = $tr.Characters($startPos + "classification ".Length, $afterLen)
    if ($afterRun.Text -eq " on a ") {
        $afterRun.Text = "on a "
    }
}
